$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D='30.048.31'; E='  -0.30%  '},
    @{Row=3; D='1.883.66'; E='  +0.33%  '},
    @{Row=4; D='0.9987'; E='  -0.06%  '},
    @{Row=5; D='243.62'; E='  -2.42%  '},
    @{Row=6; D='0.9983'; E='  -0.11%  '},
    @{Row=7; D='0.4975'; E='  -0.15%  '},
    @{Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2906'; E='  +1.96%  '},
    @{Row=9; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.06612'; E='  +0.50%  '},
    @{Row=10; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.878.95'; E='  +0.55%  '},
    @{Row=11; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='16.82'; E='  -1.84%  '},
    @{Row=12; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07199'; E='  +0.19%  '},
    @{Row=13; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6640'},
    @{Row=14; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='85.68'; E='  +0.25%  '},
    @{Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.842'; E='  +0.42%  '},
    @{Row=16; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='30.035.36'; E='  -0.27%  '},
    @{Row=17; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007766'; E='  +3.01%  '},
    @{Row=18; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.9976'; E='  -0.31%  '},
    @{Row=19; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.77'; E='  -1.28%  '},
    @{Row=20; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.120.15'; E='  +0.47%  '},
    @{Row=21; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='0.9975'; E='  -0.10%  '},
    @{Row=22; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='4.759'; E='  +0.08%  '},
    @{Row=23; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='5.593'; E='  +1.33%  '},
    @{Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='9.139'; E='  +0.89%  '},
    @{Row=25; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='149.93'; E='  +3.65%  '},
    @{Row=26; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='135.82'; E='  +0.63%  '},
    @{Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='16.76'; E='  +0.09%  '},
    @{Row=28; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.906'; E='  -2.98%  '},
    @{Row=29; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.383'; E='  -1.65%  '},
    @{Row=30; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.164'; E='  -1.54%  '},
    @{Row=31; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.08676'; E='  +0.80%  '},
    @{Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.945'; E='  +1.17%  '},
    @{Row=33; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.04996'; E='  -1.58%  '},
    @{Row=34; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.102'; E='  -3.24%  '},
    @{Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7003'; E='  +1.95%  '},
    @{Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.654'; E='  -1.64%  '},
    @{Row=37; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.691'; E='  -1.98%  '},
    @{Row=38; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.197'; E='  -6.40%  '},
    @{Row=39; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.9363'; E='  -2.87%  '},
    @{Row=40; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01643'; E='  +0.69%  '},
    @{Row=41; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.956'; E='  -2.21%  '},
    @{Row=42; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9989'; E='  -0.09%  '},
    @{Row=43; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4183'; E='  -0.33%  '},
    @{Row=44; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='101.17'; E='  -2.05%  '},
    @{Row=45; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.498'; E='  -0.56%  '},
    @{Row=46; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1259'; E='  +0.11%  '},
    @{Row=47; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05721'; E='  +1.55%  '},
    @{Row=48; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='32.31'; E='  -0.87%  '},
    @{Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='8.198'; E='  -0.70%  '},
    @{Row=50; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='55.85'; E='  +1.68%  '},
    @{Row=51; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.3706'; E='  -0.72%  '}
)

foreach ($item in $rows) {
    $r = $item.Row

    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    if ($item.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}

Write-Host "Applied $($rows.Count) row updates"
